$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UCT1")
$ws.Activate()

$win = $excel.ActiveWindow
$win.TopLeftCell = $ws.Range("Q13")
